{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies three changes to the letter body:\n//   1. \"September 19, 2025\" -> \"September 21, 2025\"\n//   2. Split the sender-address paragraph \"959 Story Road, San Jose CA 95122\"\n//      (the one in the free-flowing body text, NOT the one inside the\n//      ASSOCIATION/ACCOUNT NUMBER/PROPERTY ADDRESS table) into two\n//      paragraphs: \"959 Story Road\" and \"San Jose, CA 95122\".\n//   3. Remove the empty \"No Spacing\" paragraph directly following the\n//      \"...Board of Directors\" signature paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/tableNestingLevel,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Change 1: update the letter date -----------------------------------\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"September 19, 2025\") {\n    items[i].getRange().insertText(\"September 21, 2025\", \"Replace\");\n    break;\n  }\n}\n\n// --- Change 2: split \"959 Story Road, San Jose CA 95122\" ----------------\n// Only the body-level paragraph (outside of the table) is targeted.\nfor (let i = 0; i < items.length; i++) {\n  if (\n    items[i].text === \"959 Story Road, San Jose CA 95122\" &&\n    items[i].tableNestingLevel === 0\n  ) {\n    const para = items[i];\n    // Replace the run text with just the street line...\n    para.getRange().insertText(\"959 Story Road\", \"Replace\");\n    // ...then add a new paragraph right after it with the city/state/zip,\n    // matching the original paragraph's formatting (Arial 11pt, same as\n    // every other address line in the letterhead).\n    const newPara = para.insertParagraph(\"San Jose, CA 95122\", \"After\");\n    newPara.font.name = \"Arial\";\n    newPara.font.size = 11;\n    break;\n  }\n}\n\nawait context.sync();\n\n// --- Change 3: remove the empty paragraph after \"Board of Directors\" ----\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst items2 = paragraphs2.items;\nfor (let i = 0; i < items2.length; i++) {\n  if (items2[i].text.indexOf(\"Board of Directors\") !== -1) {\n    const next = items2[i].getNext();\n    next.load(\"text\");\n    await context.sync();\n    if (next.text === \"\") {\n      next.delete();\n    }\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies three changes to the letter body:\n#   1. \"September 19, 2025\" -> \"September 21, 2025\"\n#   2. Split the sender-address paragraph \"959 Story Road, San Jose CA 95122\"\n#      (the one in the free-flowing body text, NOT the one inside the\n#      ASSOCIATION/ACCOUNT NUMBER/PROPERTY ADDRESS table) into two\n#      paragraphs: \"959 Story Road\" and \"San Jose, CA 95122\".\n#   3. Remove the empty \"No Spacing\" paragraph directly following the\n#      \"...Board of Directors\" signature paragraph.\n\n$d = $word.ActiveDocument\n\n# Collect the character ranges covered by every table in the document, so we\n# can tell apart the free-standing body paragraph from the one living inside\n# the ASSOCIATION/ACCOUNT NUMBER/PROPERTY ADDRESS table (both paragraphs\n# otherwise contain the exact same text).\n$tableStarts = @()\n$tableEnds = @()\nfor ($t = 1; $t -le $d.Tables.Count; $t++) {\n    $tableStarts += $d.Tables($t).Range.Start\n    $tableEnds += $d.Tables($t).Range.End\n}\n\nfunction InTable($start) {\n    for ($k = 0; $k -lt $tableStarts.Count; $k++) {\n        if ($start -ge $tableStarts[$k] -and $start -lt $tableEnds[$k]) {\n            return $true\n        }\n    }\n    return $false\n}\n\n# --- Change 1: update the letter date ------------------------------------\n# NOTE: Word COM paragraph ranges include the trailing paragraph mark\n# (\"\\r\"), so text comparisons trim it off first.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text.TrimEnd() -eq \"September 19, 2025\") {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# --- Change 2: split \"959 Story Road, San Jose CA 95122\" ------------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text.TrimEnd() -eq \"959 Story Road, San Jose CA 95122\" -and -not (InTable($p.Range.Start))) {\n        # Replacing the text with the street line plus a paragraph break\n        # splits this paragraph in two, preserving the original run/\n        # paragraph formatting (Arial 11pt) on both halves.\n        $p.Range.Text = \"959 Story Road`r\"\n        $next = $d.Paragraphs($i + 1)\n        $next.Range.Text = \"San Jose, CA 95122\"\n        break\n    }\n}\n\n# --- Change 3: remove the empty paragraph after \"Board of Directors\" -----\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text.TrimEnd() -like \"*Board of Directors\") {\n        $next = $d.Paragraphs($i + 1)\n        if ($next.Range.Text.TrimEnd() -eq \"\") {\n            $next.Range.Delete()\n        }\n        break\n    }\n}\n"}
